$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D retains text storage (values like "0.524" would otherwise be
# auto-converted to numbers by Excel). Apply Text format before assigning values.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.299.23'
$ws.Range("E2").Value = '  +1.38%  '
$ws.Range("D3").Value = '1.680.50'
$ws.Range("E3").Value = '  +3.29%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = '220.86'
$ws.Range("E5").Value = '  +3.17%  '
$ws.Range("D6").Value = '0.524'
$ws.Range("E6").Value = '  +0.89%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").Value = '29.88'
$ws.Range("E8").Value = '  +0.48%  '
$ws.Range("D9").Value = '0.265'
$ws.Range("E9").Value = '  +2.10%  '
$ws.Range("E10").Value = '  +1.48%  '
$ws.Range("D11").Value = '0.0901'
$ws.Range("E11").Value = '  -1.63%  '
$ws.Range("D12").Value = '1.920.18'
$ws.Range("D13").Value = '10.64'
$ws.Range("E13").Value = '  +13.69%  '
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").Value = '0.618'
$ws.Range("E14").Value = '  +8.36%  '
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '1.673.34'
$ws.Range("E15").Value = '  +2.75%  '
$ws.Range("D16").Value = '4.00'
$ws.Range("E16").Value = '  +3.12%  '
$ws.Range("D17").Value = '30.303.89'
$ws.Range("E17").Value = '  +1.35%  '
$ws.Range("D18").Value = '65.76'
$ws.Range("E18").Value = '  +1.22%  '
$ws.Range("D19").Value = '247.36'
$ws.Range("E19").Value = '  -0.31%  '
$ws.Range("D20").Value = '0.0₃0718'
$ws.Range("E20").Value = '  +1.96%  '
$ws.Range("E21").Value = '  -0.19%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '4.30'
$ws.Range("E22").Value = '  +3.62%  '
$ws.Range("B23").Value = 'Avalanche'
$ws.Range("C23").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D23").Value = '10.20'
$ws.Range("E23").Value = '  +6.08%  '
$ws.Range("E24").Value = '  +4.45%  '
$ws.Range("D25").Value = '158.84'
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("D26").Value = '15.86'
$ws.Range("E26").Value = '  +0.99%  '
$ws.Range("E27").Value = '  +0.12%  '
$ws.Range("E28").Value = '  +2.29%  '
$ws.Range("E29").Value = '  -0.20%  '
$ws.Range("E30").Value = '  +2.02%  '
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value = '3.50'
$ws.Range("E31").Value = '  +3.89%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '1.14'
$ws.Range("E32").Value = '  +0.75%  '
$ws.Range("E33").Value = '  +3.36%  '
$ws.Range("D34").Value = '1.497.70'
$ws.Range("E34").Value = '  +4.74%  '
$ws.Range("E35").Value = '  +5.08%  '
$ws.Range("E36").Value = '  -0.35%  '
$ws.Range("E37").Value = '  +5.32%  '
$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").Value = '2.74'
$ws.Range("E38").Value = '  -4.49%  '
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").Value = '0.585'
$ws.Range("E39").Value = '  +5.38%  '
$ws.Range("B40").Value = 'Aave'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D40").Value = '78.78'
$ws.Range("E40").Value = '  +10.62%  '
$ws.Range("E41").Value = '  +1.18%  '
$ws.Range("D42").Value = '0.853'
$ws.Range("E42").Value = '  +2.77%  '
$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D43").Value = '0.0507'
$ws.Range("E43").Value = '  +2.52%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = '2.01'
$ws.Range("E44").Value = '  +1.97%  '
$ws.Range("D45").Value = '0.999'
$ws.Range("E45").Value = '  -0.09%  '
$ws.Range("E46").Value = '  -4.28%  '
$ws.Range("D47").Value = '51.86'
$ws.Range("E47").Value = '  -6.41%  '
$ws.Range("D48").Value = '1.812.96'
$ws.Range("E48").Value = '  +2.59%  '
$ws.Range("D49").Value = '5.43'
$ws.Range("E49").Value = '  -0.36%  '
$ws.Range("D50").Value = '95.27'
$ws.Range("E50").Value = '  +6.06%  '
$ws.Range("E51").Value = '  +8.45%  '
